$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 265, shifting existing rows 265-315 down to 266-316
$ws.Rows.Item(265).Insert()

# Populate the newly inserted row 265 with its data
$ws.Range("A265").Value = 5
$ws.Range("B265").Value = "Macroferia Regional de Talca"
$ws.Range("C265").Value = "Maule"
$ws.Range("D265").Value = 45211
$ws.Range("E265").Value = 7
$ws.Range("F265").Value = 100112028
$ws.Range("G265").Value = "Sandia"
$ws.Range("H265").Value = "Sin especificar"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 3000
$ws.Range("K265").Value = 600
$ws.Range("L265").Value = 600
$ws.Range("M265").Value = 600
$ws.Range("N265").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O265").Value = "Perú"
$ws.Range("P265").Value = 600
$ws.Range("Q265").Value = 1
$ws.Range("R265").Value = "Hortaliza"
